# edit.ps1 - applies the "Update gh-pages to output generated at 456a3b4" diff
# to 上海-漫展信息.xlsx (4 sheets: 展览 / 演出 / 本地生活 / 全部类型).
#
# Most changes are simple numeric updates to column F ("想去人数" / interest
# counter). In the "全部类型" sheet, a new concert entry was also inserted
# at row 34 (pushing several rows down by one) while the old row-40 entry
# ("坏孩纸物语...") was removed, so rows 34-41 there get their full row
# content rewritten instead of just column F.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws = $wb.Worksheets("展览")
$ws.Range("F3").Value = 260
$ws.Range("F4").Value = 584
$ws.Range("F5").Value = 2559
$ws.Range("F6").Value = 10
$ws.Range("F9").Value = 247
$ws.Range("F10").Value = 5237
$ws.Range("F11").Value = 97
$ws.Range("F12").Value = 1452
$ws.Range("F13").Value = 1376
$ws.Range("F15").Value = 6940
$ws.Range("F16").Value = 387
$ws.Range("F17").Value = 46
$ws.Range("F20").Value = 4658
$ws.Range("F22").Value = 69
$ws.Range("F23").Value = 2308
$ws.Range("F24").Value = 1248
$ws.Range("F25").Value = 436
$ws.Range("F26").Value = 1150
$ws.Range("F27").Value = 214
$ws.Range("F28").Value = 90
$ws.Range("F29").Value = 69
$ws.Range("F30").Value = 155
$ws.Range("F32").Value = 1266
$ws.Range("F33").Value = 1980
$ws.Range("F34").Value = 224
$ws.Range("F35").Value = 512
$ws.Range("F36").Value = 199
$ws.Range("F37").Value = 1364
$ws.Range("F39").Value = 86
$ws.Range("F40").Value = 101
$ws.Range("F41").Value = 160
$ws.Range("F42").Value = 1109
$ws.Range("F43").Value = 2398
$ws.Range("F45").Value = 64
$ws.Range("F47").Value = 228
$ws.Range("F48").Value = 69
$ws.Range("F49").Value = 11

# --- Sheet: 演出 ---
$ws = $wb.Worksheets("演出")
$ws.Range("F5").Value = 452
$ws.Range("F7").Value = 133
$ws.Range("F11").Value = 4
$ws.Range("F12").Value = 380
$ws.Range("F13").Value = 266
$ws.Range("F16").Value = 180
$ws.Range("F25").Value = 2
$ws.Range("F28").Value = 277
$ws.Range("F29").Value = 25

# --- Sheet: 本地生活 ---
$ws = $wb.Worksheets("本地生活")
$ws.Range("F6").Value = 1650
$ws.Range("F8").Value = 1277
$ws.Range("F10").Value = 1732
$ws.Range("F11").Value = 2168
$ws.Range("F12").Value = 592
$ws.Range("F13").Value = 496

# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets("全部类型")
$ws.Range("F3").Value = 1650
$ws.Range("F4").Value = 260
$ws.Range("F5").Value = 584
$ws.Range("F7").Value = 2559
$ws.Range("F9").Value = 1277
$ws.Range("F10").Value = 2168
$ws.Range("F11").Value = 5238
$ws.Range("F12").Value = 592
$ws.Range("F13").Value = 452
$ws.Range("F14").Value = 133
$ws.Range("F15").Value = 97
$ws.Range("F17").Value = 1452
$ws.Range("F18").Value = 1376
$ws.Range("F20").Value = 6940
$ws.Range("F21").Value = 387
$ws.Range("F22").Value = 497
$ws.Range("F23").Value = 46
$ws.Range("F24").Value = 4658
$ws.Range("F25").Value = 2308
$ws.Range("F26").Value = 1248
$ws.Range("F27").Value = 436
$ws.Range("F28").Value = 1150
$ws.Range("F29").Value = 214
$ws.Range("F30").Value = 69
$ws.Range("F31").Value = 266
$ws.Range("F33").Value = 155
$ws.Range("F42").Value = 160
$ws.Range("F43").Value = 25
$ws.Range("F44").Value = 1109
$ws.Range("F45").Value = 2398
$ws.Range("F46").Value = 64
$ws.Range("F47").Value = 228
$ws.Range("F48").Value = 69

# Row-shift region: new row inserted at 34 ("8.17 跨越国境的旋律..." concert),
# rows 34-37 (LOVELIVE / BACG / ICOMIC / 魔都特摄) shift down to 35-38,
# row 38 (HAG concert) shifts to 39, row 39 (coser) shifts to 40,
# old row 40 ("坏孩纸物语...") is dropped, row 41 (wanuka) stays put (F only).
# Row 34: C, D, E, F, G, H, I changed
$ws.Range("C34").Value = "上海·8.17 跨越国境的旋律——【吹响吧！上低音号】音乐监修、洗足音大教授大和田雅洋x和音社交响吹奏音乐会"
$ws.Range("D34").Value = "复兴中路1380号 捷豹上海交响音乐厅"
$ws.Range("E34").Value = "2024.08.17 19:30-08.17 21:15"
$ws.Range("F34").Value = 180
$ws.Range("G34").Value = 80
$ws.Range("H34").Value = "https://show.bilibili.com/platform/detail.html?id=89733"
$ws.Range("I34").Value = "//i0.hdslb.com/bfs/openplatform/202407/MeZ4tVLG1721651973268.jpeg"

# Row 35: C, D, E, F, G, H, I changed
$ws.Range("C35").Value = "上海·LOVELIVE 同人ONLY"
$ws.Range("D35").Value = "海潮路133号B1 JUMP工坊"
$ws.Range("E35").Value = "2024.08.17 14:00-08.17 19:00"
$ws.Range("F35").Value = 364
$ws.Range("G35").Value = 60
$ws.Range("H35").Value = "https://show.bilibili.com/platform/detail.html?id=86711"
$ws.Range("I35").Value = "//i2.hdslb.com/bfs/openplatform/202405/bllJHQFL1716983812432.jpeg"

# Row 36: C, D, E, F, G, H, I changed
$ws.Range("C36").Value = "上海·第六届燃梦BACG PRO动漫嘉年华-我们在燃梦相遇吧！"
$ws.Range("D36").Value = "盈浦街道淀山浦社区淀山湖大道851号青浦万达茂F3 万达汽车乐园(青浦万达茂店)"
$ws.Range("E36").Value = "2024.08.17 11:00-08.18 16:00"
$ws.Range("F36").Value = 1980
$ws.Range("G36").Value = 65.8
$ws.Range("H36").Value = "https://show.bilibili.com/platform/detail.html?id=85239"
$ws.Range("I36").Value = "//i1.hdslb.com/bfs/openplatform/202405/mzD4rhY21715109458100.jpeg"

# Row 37: C, D, E, F, G, H, I changed
$ws.Range("C37").Value = "上海·首届ICOMIC漫展·动漫游戏嘉年华.让我们追随热爱·与你相遇"
$ws.Range("D37").Value = "红宝石路188号 古北SOHO"
$ws.Range("E37").Value = "2024.08.17 10:00-08.18 17:00"
$ws.Range("F37").Value = 224
$ws.Range("G37").Value = 68
$ws.Range("H37").Value = "https://show.bilibili.com/platform/detail.html?id=88975"
$ws.Range("I37").Value = "//i0.hdslb.com/bfs/openplatform/202407/GZMx5q0L1722220472815.jpeg"

# Row 38: B, C, D, E, F, G, H, I changed
$ws.Range("B38").Value = "2024-08-17"
$ws.Range("C38").Value = "上海·魔都特摄同人嘉年华"
$ws.Range("D38").Value = "天山路1111号 现所创邑MIX"
$ws.Range("E38").Value = "2024.08.17 09:30-08.18 17:30"
$ws.Range("F38").Value = 512
$ws.Range("G38").Value = 69
$ws.Range("H38").Value = "https://show.bilibili.com/platform/detail.html?id=89516"
$ws.Range("I38").Value = "//i0.hdslb.com/bfs/openplatform/202407/0050E5641721292312668.png"

# Row 39: C, D, E, F, G, H, I changed
$ws.Range("C39").Value = "上海·HAG 1st live in Shanghai《不眨眼》2024演唱会"
$ws.Range("D39").Value = "中兴路1683号金融街购物中心三楼L3-27 蜚声LIVE House"
$ws.Range("E39").Value = "2024.08.24 19:30-08.24 21:30"
$ws.Range("F39").Value = 26
$ws.Range("G39").Value = 480
$ws.Range("H39").Value = "https://show.bilibili.com/platform/detail.html?id=89977"
$ws.Range("I39").Value = "//i1.hdslb.com/bfs/openplatform/202407/iXZNZNM01722243246403.png"

# Row 40: B, C, D, E, F, G, H, I changed
$ws.Range("B40").Value = "2024-08-24"
$ws.Range("C40").Value = "上海·coser动漫展"
$ws.Range("D40").Value = "海潮路133号B1 JUMP工坊"
$ws.Range("E40").Value = "2024.08.24 10:00-08.25 17:00"
$ws.Range("F40").Value = 1364
$ws.Range("G40").Value = 60
$ws.Range("H40").Value = "https://show.bilibili.com/platform/detail.html?id=87347"
$ws.Range("I40").Value = "//i0.hdslb.com/bfs/openplatform/202406/i6vAgX8I1719311206769.jpeg"

# Row 41: F changed
$ws.Range("F41").Value = 2

